# The "FSE Architecture" deck has a rounded-rectangle shape on slide 2 whose
# text reads "LLM Prediction" (a single bold run). The commit retitles it to
# "LLM Answer Generation", realized as two runs: the existing "LLM " text is
# left alone and the word "Prediction" is replaced with "Answer Generation".
#
# Find the shape robustly (by its current text) instead of hard-coding a
# slide/shape index, then replace just the "Prediction" substring so the
# leading "LLM " run keeps its original run properties untouched.

$p = $ppt.ActivePresentation

$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "LLM Prediction") {
                $targetShape = $shp
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange
    # "Prediction" starts right after "LLM " -> 1-based char 5, length 10.
    $old = $tr.Characters(5, 10)
    $old.Text = "Answer Generation"
}
